$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add a bottom border under the (previously last) data row ---
$ws.Range("A3:E3").Borders.Item(9).LineStyle = 1
$ws.Range("A3:E3").Borders.Item(9).Weight = 2
# A3 stays empty, but now carries the new bordered style
$ws.Range("A3").Value = ""

# --- New rows 4-6: English / Russian / "converted" triplets ---
# Row numbers (column B)
$ws.Range("B4").Value = 307
$ws.Range("B5").Value = 310
$ws.Range("B6").Value = 313

# Column C - English text (fills shared strings 12,13,14)
$ws.Range("C4").Value = " Welcome back!"
$ws.Range("C5").Value = " Oh! After that expedition, you\nboth seem a little different. My dears, you\nhave a new lean, mean look!"
$ws.Range("C6").Value = " I can see that!"

# Column D - Russian translated text (fills shared strings 15,16,17)
$ws.Range("D4").Value = " С возвращением!"
$ws.Range("D5").Value = " Ой! Вы выглядите иначе после\nэкспедиции. Дорогуши мои, вы совсем\nотощали!"
$ws.Range("D6").Value = " Это заметно! "

# Column E - converted/obfuscated text (fills shared strings 18,19,20)
$ws.Range("E4").Value = " Ò âïèâñàþåîéåí!"
$ws.Range("E5").Value = " Ïê! Âú âúãìÿäéóå éîàœå ðïòìå\nüëòðåäéøéé. Äïñïãôšé íïé, âú òïâòåí\nïóïþàìé!"
$ws.Range("E6").Value = " Üóï èàíåóîï!"

# Row 5 has the taller line height seen in the diff
$ws.Rows.Item(5).RowHeight = 31.8

# Move the active selection to D5, matching the post-edit state
$ws.Range("D5").Select()
